# Apply Betfair Back/Lay odds updates for 2025-10-06 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 2.56
$ws.Range("J2").Value = 3.9
$ws.Range("L2").Value = 1.26
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 4.6
$ws.Range("O2").Value = 1.18
$ws.Range("P2").Value = 2.56
$ws.Range("Q2").Value = 1.52
$ws.Range("R2").Value = 1.61
$ws.Range("S2").Value = 2.4
$ws.Range("T2").Value = 1.03
$ws.Range("U2").Value = 2.62
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 17
$ws.Range("Z2").Value = 22
$ws.Range("AB2").Value = 980
$ws.Range("AC2").Value = 10.5
$ws.Range("AH2").Value = 15.5
$ws.Range("AM2").Value = 60
$ws.Range("AO2").Value = 14.5

# Row 3
$ws.Range("T3").Value = 2.04

# Row 5
$ws.Range("F5").Value = 1.33
$ws.Range("G5").Value = 1.48
$ws.Range("H5").Value = 1.09
$ws.Range("I5").Value = 13.5
$ws.Range("K5").Value = 6.4
$ws.Range("S5").Value = 1.6
$ws.Range("V5").Value = 1.08
$ws.Range("W5").Value = 3.05

# Row 8
$ws.Range("N8").Value = 1.3
$ws.Range("P8").Value = 1.3
$ws.Range("W8").Value = 1.02

# Row 9
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 2.52
$ws.Range("I9").Value = 2.7
$ws.Range("J9").Value = 3.7
$ws.Range("K9").Value = 3.8
$ws.Range("P9").Value = 1.35
$ws.Range("V9").Value = 1.58
$ws.Range("X9").Value = 15.5
$ws.Range("Y9").Value = 13
$ws.Range("Z9").Value = 18
$ws.Range("AA9").Value = 40
$ws.Range("AC9").Value = 8.199999999999999
$ws.Range("AD9").Value = 12.5
$ws.Range("AE9").Value = 27
$ws.Range("AF9").Value = 21
$ws.Range("AH9").Value = 16.5
$ws.Range("AI9").Value = 40
$ws.Range("AJ9").Value = 46
$ws.Range("AK9").Value = 30
$ws.Range("AL9").Value = 44
$ws.Range("AM9").Value = 80
$ws.Range("AN9").Value = 25
$ws.Range("AO9").Value = 22

# Row 10
$ws.Range("G10").Value = 2.96
$ws.Range("H10").Value = 2.64
$ws.Range("K10").Value = 3.7
$ws.Range("Z10").Value = 980
$ws.Range("AA10").Value = 980
$ws.Range("AF10").Value = 980
$ws.Range("AH10").Value = 980
$ws.Range("AI10").Value = 980
$ws.Range("AJ10").Value = 980
$ws.Range("AK10").Value = 980
$ws.Range("AL10").Value = 980
$ws.Range("AN10").Value = 980
$ws.Range("AO10").Value = 980

# Row 11
$ws.Range("Q11").Value = 2.36
$ws.Range("T11").Value = 2.04
$ws.Range("U11").Value = 1.56
$ws.Range("X11").Value = 980
$ws.Range("Y11").Value = 980
$ws.Range("Z11").Value = 980
$ws.Range("AB11").Value = 980
$ws.Range("AC11").Value = 980
$ws.Range("AD11").Value = 980
$ws.Range("AF11").Value = 980
$ws.Range("AG11").Value = 980
$ws.Range("AO11").Value = 980

# Row 12
$ws.Range("F12").Value = 3.55
$ws.Range("G12").Value = 3.9
$ws.Range("I12").Value = 2.56
$ws.Range("K12").Value = 3.25
$ws.Range("M12").Value = 1.16
$ws.Range("S12").Value = 7.6
$ws.Range("T12").Value = 2.52
$ws.Range("W12").Value = 1.34
$ws.Range("Z12").Value = 980
$ws.Range("AA12").Value = 980
$ws.Range("AB12").Value = 8.4
$ws.Range("AC12").Value = 980
$ws.Range("AD12").Value = 980
$ws.Range("AE12").Value = 980
$ws.Range("AF12").Value = 980
$ws.Range("AG12").Value = 980
$ws.Range("AH12").Value = 980
$ws.Range("AJ12").Value = 120
$ws.Range("AK12").Value = 80
$ws.Range("AL12").Value = 150
$ws.Range("AN12").Value = 170

# Row 13
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 2.22
$ws.Range("H13").Value = 4.3
$ws.Range("I13").Value = 5.1
$ws.Range("J13").Value = 3
$ws.Range("K13").Value = 3.75
$ws.Range("M13").Value = 1.12
$ws.Range("N13").Value = 2.54
$ws.Range("O13").Value = 1.53
$ws.Range("P13").Value = 1.51
$ws.Range("Q13").Value = 2.58
$ws.Range("R13").Value = 1.18
$ws.Range("S13").Value = 5.3
$ws.Range("T13").Value = 2.14
$ws.Range("U13").Value = 1.71
$ws.Range("V13").Value = 1.24
$ws.Range("W13").Value = 1.81
$ws.Range("X13").Value = 980

# Row 14
$ws.Range("H14").Value = 1.04
$ws.Range("J14").Value = 3.55
$ws.Range("L14").Value = 1.01
$ws.Range("M14").Value = 1.01
$ws.Range("N14").Value = 1.11
$ws.Range("O14").Value = 1.01
$ws.Range("R14").Value = 1.14
$ws.Range("S14").Value = 1.92
$ws.Range("T14").Value = 1.03
$ws.Range("U14").Value = 1.03
$ws.Range("V14").Value = 1.02
$ws.Range("W14").Value = 2.38
$ws.Range("X14").Value = 1000
$ws.Range("Y14").Value = 1000
$ws.Range("Z14").Value = 1000
$ws.Range("AA14").Value = 1000
$ws.Range("AB14").Value = 1000
$ws.Range("AC14").Value = 1000
$ws.Range("AD14").Value = 1000
$ws.Range("AE14").Value = 1000
$ws.Range("AF14").Value = 1000
$ws.Range("AG14").Value = 1000
$ws.Range("AH14").Value = 1000
$ws.Range("AI14").Value = 1000
$ws.Range("AJ14").Value = 1000
$ws.Range("AK14").Value = 1000
$ws.Range("AL14").Value = 1000
$ws.Range("AM14").Value = 1000
$ws.Range("AN14").Value = 1000
$ws.Range("AO14").Value = 1000
